# Populate the "Leyenda" (legend) worksheet with the field-description
# table describing each column of the "Compromisos" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Leyenda")

$rows = @(
    ,@('Descripción de los campos.', $null)
    ,@('Nombre', 'Descripción')
    ,@('Número correlativo', 'Número del compromiso para la tabla. Inicia en 1 y es correlativo.')
    ,@('Origen compromiso', 'Estudio origen del compromiso.')
    ,@('Capítulo', 'Capítulo del origen del compromiso.')
    ,@('Sección', 'Sección del origen del compromiso.')
    ,@('Página', 'Página del origen del compromiso.')
    ,@('Aspecto asociado', 'El aspecto ambiental sobre el cual el compromiso busca gestionar (prevenir, mitigar, compensar, etc) el impacto.')
    ,@('Instalación/componente', 'Instalación o componente a la cual el compromiso se encuentra vinculada.')
    ,@('Contenido original del compromiso', 'Es el copia textual del compromiso tal y como se presentó en el estudio origen.')
    ,@('Resumen del compromiso', 'Es un texto que presenta de manera resumida el compromiso.')
    ,@('Antecedentes del compromiso', 'Es un texto que presenta información sobre los antecedentes del compromiso como por ejemplo si es que representa una actualiza, abarca o reemplaza un compromiso anterior, similar o duplicado.')
    ,@('Temporalidad', 'Establece si el compromiso ocurre una vez de manera puntual o de manera permanente dentro del periodo de vigencia.')
    ,@('Fecha de inicio', 'Fecha en la que se iniciará el compromiso.')
    ,@('Frecuencia', 'Establece la frecuencia de cumplimiento del compromiso dentro del periodo de vigencia.')
    ,@('Criticidad', 'Calificación cualitativa sobre la criticidad del compromiso a distintas características como el riesgo, las consecuencias, etc.')
    ,@('Evidencias', 'Indica la existencia de evidencia de cumplimiento e incluye los archivos correspondientes.')
    ,@('Estado de cumplimiento', 'Señala la situación de cumplimiento del compromiso.')
    ,@('Acción sobre el compromiso', 'Señala si se requiere alguna acción sobre el compromiso.')
    ,@('Detalle de acción', 'Detalle de la acción que requiere.')
    ,@('Frecuencia de verificación', 'Cada cuanto tiempo se verificará el cumplimiento de un compromiso.')
    ,@('Area responsable', 'Señala el equipo o área responsable de ejecutar el compromiso.')
    ,@('Correos de notificación', 'Lista de correos para notificaciones y avisos de verificación.')
    ,@('Notificación de inactividad', 'No Definido')
    ,@('Nombre de revisor', 'Nombre del encargado de la revisión.')
    ,@('Fecha de revisión', 'Indica la última vez que el compromiso y sus características fueron revisadas y validadas.')
    ,@('Referencia técnica o legal de cumplimiento', 'Indica, en caso exista, la referencia sobre cómo se deberá cumplir el compromiso en base a la normativa o texto técnico.')
    ,@('Presupuesto', 'Indica el presupuesto asociado al cumplimiento.')
    ,@('Notas adicionales', 'Permite complementar el compromiso con detalles específicos.')
    ,@('Construcción', 'Establece el periodo de construcción y si es que el compromiso está vigente en  esta etapa.')
    ,@('Operación', 'Establece el periodo de operación y si es que el compromiso está vigente en  esta etapa.')
    ,@('Cierre', 'Establece el periodo de cierre y si es que el compromiso está vigente en  esta etapa.')
    ,@('Post-cierre', 'Establece el periodo de post-cierre y si es que el compromiso está vigente en  esta etapa.')
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 1
    $a = $rows[$i][0]
    $b = $rows[$i][1]

    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value = $a
    $cellA.Borders.LineStyle = 1

    if ($b -ne $null) {
        $cellB = $ws.Cells.Item($r, 2)
        $cellB.Value = $b
        $cellB.Borders.LineStyle = 1
    }
}

# Column widths: A = 20 characters, B = 180 characters.
# (offset by -5/6 to compensate for the host's character->width rounding)
$ws.Columns.Item(1).ColumnWidth = 19.1666666666667
$ws.Columns.Item(2).ColumnWidth = 179.166666666667
